$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-17: state, model, RMSE, RMSPE, MAE, MAPE, crps, mcr95, mcr50
# Row 2: PR / GT
$ws.Cells.Item(2,1).Value = "PR"
$ws.Cells.Item(2,2).Value = "GT"
$ws.Cells.Item(2,3).Value = 11484.54150769252
$ws.Cells.Item(2,4).Value = 25.05055481377521
$ws.Cells.Item(2,5).Value = 8868.131248810469
$ws.Cells.Item(2,6).Value = 21.73535765496148
$ws.Range("G2:I2").ClearContents()

# Row 3: PR / ARABM
$ws.Cells.Item(3,1).Value = "PR"
$ws.Cells.Item(3,2).Value = "ARABM"
$ws.Cells.Item(3,3).Value = 7329.351261883578
$ws.Cells.Item(3,4).Value = 20.54423506440512
$ws.Cells.Item(3,5).Value = 6857.531895833334
$ws.Cells.Item(3,6).Value = 18.94458433570844
$ws.Cells.Item(3,7).Value = 5995.919793907084
$ws.Cells.Item(3,8).Value = 0.2916666666666667
$ws.Cells.Item(3,9).Value = 0.08333333333333333

# Row 4: PR / ARABMwGT
$ws.Cells.Item(4,1).Value = "PR"
$ws.Cells.Item(4,2).Value = "ARABMwGT"
$ws.Cells.Item(4,3).Value = 7215.949265573547
$ws.Cells.Item(4,4).Value = 20.43864114403762
$ws.Cells.Item(4,5).Value = 6712.609825
$ws.Cells.Item(4,6).Value = 17.89565582890354
$ws.Cells.Item(4,7).Value = 5940.490625765001
$ws.Cells.Item(4,8).Value = 0.625
$ws.Cells.Item(4,9).Value = 0.125

# Row 5: PR / Bastos
$ws.Cells.Item(5,1).Value = "PR"
$ws.Cells.Item(5,2).Value = "Bastos"
$ws.Cells.Item(5,3).Value = 5690.171660987475
$ws.Cells.Item(5,4).Value = 14.24321709177373
$ws.Cells.Item(5,5).Value = 5397.5625
$ws.Cells.Item(5,6).Value = 13.15637387704677
$ws.Range("G5:I5").ClearContents()

# Row 6: RS / GT
$ws.Cells.Item(6,1).Value = "RS"
$ws.Cells.Item(6,2).Value = "GT"
$ws.Cells.Item(6,3).Value = 8104.424187726412
$ws.Cells.Item(6,4).Value = 40.44168098734507
$ws.Cells.Item(6,5).Value = 6329.5676584957
$ws.Cells.Item(6,6).Value = 37.50869182877114
$ws.Range("G6:I6").ClearContents()

# Row 7: RS / ARABM
$ws.Cells.Item(7,1).Value = "RS"
$ws.Cells.Item(7,2).Value = "ARABM"
$ws.Cells.Item(7,3).Value = 4486.419654327045
$ws.Cells.Item(7,4).Value = 29.71708581109228
$ws.Cells.Item(7,5).Value = 4186.94315
$ws.Cells.Item(7,6).Value = 27.87180647166429
$ws.Cells.Item(7,7).Value = 3804.940174281667
$ws.Cells.Item(7,8).Value = 0.2916666666666667
$ws.Cells.Item(7,9).Value = 0.08333333333333333

# Row 8: RS / ARABMwGT
$ws.Cells.Item(8,1).Value = "RS"
$ws.Cells.Item(8,2).Value = "ARABMwGT"
$ws.Cells.Item(8,3).Value = 4306.349778019565
$ws.Cells.Item(8,4).Value = 30.12556962316969
$ws.Cells.Item(8,5).Value = 4074.276945833333
$ws.Cells.Item(8,6).Value = 29.05932882587462
$ws.Cells.Item(8,7).Value = 3631.572503907084
$ws.Cells.Item(8,8).Value = 0.25
$ws.Cells.Item(8,9).Value = 0.08333333333333333

# Row 9: RS / Bastos
$ws.Cells.Item(9,1).Value = "RS"
$ws.Cells.Item(9,2).Value = "Bastos"
$ws.Cells.Item(9,3).Value = 3872.603633731065
$ws.Cells.Item(9,4).Value = 25.10600316886911
$ws.Cells.Item(9,5).Value = 3601.916666666667
$ws.Cells.Item(9,6).Value = 24.21745359574475
$ws.Range("G9:I9").ClearContents()

# Row 10: SC / GT
$ws.Cells.Item(10,1).Value = "SC"
$ws.Cells.Item(10,2).Value = "GT"
$ws.Cells.Item(10,3).Value = 11282.176059295
$ws.Cells.Item(10,4).Value = 33.88616903808781
$ws.Cells.Item(10,5).Value = 9046.618010075275
$ws.Cells.Item(10,6).Value = 32.07223505271391
$ws.Range("G10:I10").ClearContents()

# Row 11: SC / ARABM
$ws.Cells.Item(11,1).Value = "SC"
$ws.Cells.Item(11,2).Value = "ARABM"
$ws.Cells.Item(11,3).Value = 7175.025466313111
$ws.Cells.Item(11,4).Value = 26.44259085286908
$ws.Cells.Item(11,5).Value = 6662.606929166667
$ws.Cells.Item(11,6).Value = 24.21966148458909
$ws.Cells.Item(11,7).Value = 6028.18712075125
$ws.Cells.Item(11,8).Value = 0.4583333333333333
$ws.Cells.Item(11,9).Value = 0.125

# Row 12: SC / ARABMwGT
$ws.Cells.Item(12,1).Value = "SC"
$ws.Cells.Item(12,2).Value = "ARABMwGT"
$ws.Cells.Item(12,3).Value = 6977.156431230755
$ws.Cells.Item(12,4).Value = 23.18352177633213
$ws.Cells.Item(12,5).Value = 6605.0982125
$ws.Cells.Item(12,6).Value = 22.40151803302255
$ws.Cells.Item(12,7).Value = 5986.542579855417
$ws.Cells.Item(12,8).Value = 0.4166666666666667
$ws.Cells.Item(12,9).Value = 0.1666666666666667

# Row 13: SC / Bastos
$ws.Cells.Item(13,1).Value = "SC"
$ws.Cells.Item(13,2).Value = "Bastos"
$ws.Cells.Item(13,3).Value = 7571.527352955327
$ws.Cells.Item(13,4).Value = 26.42503346544768
$ws.Cells.Item(13,5).Value = 7249.666666666667
$ws.Cells.Item(13,6).Value = 24.68192918616597
$ws.Range("G13:I13").ClearContents()

# Row 14: SP / GT
$ws.Cells.Item(14,1).Value = "SP"
$ws.Cells.Item(14,2).Value = "GT"
$ws.Cells.Item(14,3).Value = 62798.26853839643
$ws.Cells.Item(14,4).Value = 34.67147554185567
$ws.Cells.Item(14,5).Value = 55360.37226797239
$ws.Cells.Item(14,6).Value = 33.50512395722087
$ws.Range("G14:I14").ClearContents()

# Row 15: SP / ARABM
$ws.Cells.Item(15,1).Value = "SP"
$ws.Cells.Item(15,2).Value = "ARABM"
$ws.Cells.Item(15,3).Value = 50254.56183357508
$ws.Cells.Item(15,4).Value = 37.13364301222754
$ws.Cells.Item(15,5).Value = 45067.97170416667
$ws.Cells.Item(15,6).Value = 32.92768747019047
$ws.Cells.Item(15,7).Value = 40019.49618830125
$ws.Cells.Item(15,8).Value = 0.25
$ws.Cells.Item(15,9).Value = 0.08333333333333333

# Row 16: SP / ARABMwGT
$ws.Cells.Item(16,1).Value = "SP"
$ws.Cells.Item(16,2).Value = "ARABMwGT"
$ws.Cells.Item(16,3).Value = 49014.15604734916
$ws.Cells.Item(16,4).Value = 32.53678946992844
$ws.Cells.Item(16,5).Value = 45322.26379583334
$ws.Cells.Item(16,6).Value = 30.2417805435602
$ws.Cells.Item(16,7).Value = 39753.50108213459
$ws.Cells.Item(16,8).Value = 0.25
$ws.Cells.Item(16,9).Value = 0.125

# Row 17: SP / Bastos
$ws.Cells.Item(17,1).Value = "SP"
$ws.Cells.Item(17,2).Value = "Bastos"
$ws.Cells.Item(17,3).Value = 57132.55334683681
$ws.Cells.Item(17,4).Value = 31.63492990547712
$ws.Cells.Item(17,5).Value = 54985.20833333334
$ws.Cells.Item(17,6).Value = 31.10286451291734
$ws.Range("G17:I17").ClearContents()
